# Test script placeholder
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
Write-Output $s.Shapes.Count
